$d = $word.ActiveDocument

# 1. Title / heading text - remove " - Review" suffix (appears twice: Heading1 and bold run near end)
$d.Content.Find.Execute("Play Montezuma Megaways Buy Pass for Free - Review", $true, $false, $false, $false, $false, $true, 1, $false, "Play Montezuma Megaways Buy Pass for Free", 2)

# 2. "What we like" bullet 1
$d.Content.Find.Execute("Unique gameplay mechanism", $true, $false, $false, $false, $false, $true, 1, $false, "Unique and challenging gameplay mechanism", 2)

# 3. "What we like" bullet 3
$d.Content.Find.Execute("Visually appealing with vivid colors", $true, $false, $false, $false, $false, $true, 1, $false, "Vivid colors and visually appealing graphics", 2)

# 4. "What we like" bullet 4
$d.Content.Find.Execute("Refreshing Aztec theme with well-thought-out setting", $true, $false, $false, $false, $false, $true, 1, $false, "Refreshing Aztec theme and setting", 2)

# 5. "What we don't like" bullet 1
$d.Content.Find.Execute("Lacks variety in bonus features compared to other slot games", $true, $false, $false, $false, $false, $true, 1, $false, "Limited slot game options for fans of other themes", 2)

# 6. "What we don't like" bullet 2
$d.Content.Find.Execute("Can be highly volatile, making it less accessible for casual players", $true, $false, $false, $false, $false, $true, 1, $false, "May be overwhelming for beginners due to complex mechanics", 2)

# 7. Meta description italic run at the end
$d.Content.Find.Execute("Discover the gameplay, graphics, and theme of Montezuma Megaways Buy Pass. Read our review and play for free to experience this unique Aztec slot game.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of the game Montezuma Megaways Buy Pass and play it for free.", 2)
